# Apply "indy hub changed col names again" edits to the test-case-trend sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header columns (row 1) ---
$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "m1e_covid_tests"
$ws.Range("C1").Value = "m1e_daily_delta_tests"
$ws.Range("D1").Value = "m1e_daily_base_tests"
$ws.Range("E1").Value = "m1e_covid_deaths"
$ws.Range("F1").Value = "m1e_daily_delta_deaths"
$ws.Range("G1").Value = "m1e_daily_base_deaths"
$ws.Range("H1").Value = "m1e_covid_cases"
$ws.Range("I1").Value = "m1e_daily_delta_cases"
$ws.Range("J1").Value = "m1e_daily_base_cases"
$ws.Range("K1").Value = "m1e_cumsum_covid_cases"
$ws.Range("L1").Value = "m1e_cumsum_covid_deaths"
$ws.Range("M1").Value = "m1e_cumsum_covid_tests"

# --- 2. Update numeric data for rows 94-135 and append new row 136 ---
# Each entry: row number, then values for columns B..M (A/date stays the same
# except for the brand-new row 136).
$rows = @(
    @(94,  8129,0,8129,20,0,20,411,0,411,34887,2089,299318),
    @(95,  5715,0,5715,15,0,15,477,0,477,35364,2104,305033),
    @(96,  6756,0,6756,24,0,24,426,0,426,35790,2128,311789),
    @(97,  6973,0,6973,22,0,22,501,0,501,36291,2150,318762),
    @(98,  7938,0,7938,21,0,21,420,0,420,36711,2171,326700),
    @(99,  3715,0,3715,8,0,8,402,0,402,37113,2179,330415),
    @(100, 1990,0,1990,15,0,15,233,0,233,37346,2194,332405),
    @(101, 9703,0,9703,16,0,16,386,0,386,37732,2210,342108),
    @(102, 6599,1,6598,24,0,24,306,0,306,38038,2234,348707),
    @(109, 9812,2,9810,13,0,13,237,0,237,40619,2334,404703),
    @(110, 10147,3,10144,17,0,17,550,0,550,41169,2351,414850),
    @(111, 10127,0,10127,23,0,23,368,0,368,41537,2374,424977),
    @(112, 10219,7,10212,10,0,10,381,0,381,41918,2384,435196),
    @(113, 4391,1,4390,10,0,10,371,0,371,42289,2394,439587),
    @(114, 1834,3,1831,11,0,11,248,0,248,42537,2405,441421),
    @(115, 10631,3,10628,14,0,14,304,0,304,42841,2419,452052),
    @(116, 9642,3,9639,5,0,5,273,0,273,43114,2424,461694),
    @(117, 9082,2,9080,9,0,9,511,0,511,43625,2433,470776),
    @(118, 9278,0,9278,12,0,12,445,0,445,44070,2445,480054),
    @(119, 9712,4,9708,11,0,11,478,0,478,44548,2456,489766),
    @(120, 3975,13,3962,10,0,10,348,0,348,44896,2466,493741),
    @(121, 2534,3,2531,9,0,9,300,0,300,45196,2475,496275),
    @(122, 9341,55,9286,7,0,7,371,0,371,45567,2482,505616),
    @(123, 8524,37,8487,7,0,7,357,0,357,45924,2489,514140),
    @(124, 6983,70,6913,10,0,10,440,0,440,46364,2499,521123),
    @(125, 8097,16,8081,7,0,7,528,0,528,46892,2506,529220),
    @(126, 4639,97,4542,7,1,6,497,0,497,47389,2513,533859),
    @(127, 2031,8,2023,6,0,6,589,0,589,47978,2519,535890),
    @(128, 2191,5,2186,12,0,12,315,0,315,48293,2531,538081),
    @(129, 10118,60,10058,11,0,11,301,0,301,48594,2542,548199),
    @(130, 8441,327,8114,6,1,5,445,0,445,49039,2548,556640),
    @(131, 7496,666,6830,7,0,7,499,0,499,49538,2555,564136),
    @(132, 6286,1772,4514,8,0,8,726,0,726,50264,2563,570422),
    @(133, 5323,3138,2185,5,2,3,763,0,763,51027,2568,575745),
    @(134, 1540,715,825,6,1,5,550,0,550,51577,2574,577285),
    @(135, 945,829,116,6,6,0,447,1,446,52024,2580,578230)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $ws.Cells.Item($r, 12).Value = $row[11]
    $ws.Cells.Item($r, 13).Value = $row[12]
}

# --- 3. Append brand-new row 136 ---
# Force the date-like string to stay as text (not auto-converted to a date
# serial number) and drop any stray number-format style afterwards.
$ws.Cells.Item(136, 1).NumberFormat = "@"
$ws.Cells.Item(136, 1).Value = "2020-07-13"
$ws.Cells.Item(136, 1).Style = "Normal"
$ws.Cells.Item(136, 2).Value = 179
$ws.Cells.Item(136, 3).Value = 179
$ws.Cells.Item(136, 4).Value = 0
$ws.Cells.Item(136, 5).Value = 2
$ws.Cells.Item(136, 6).Value = 2
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 661
$ws.Cells.Item(136, 9).Value = 661
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 52685
$ws.Cells.Item(136, 12).Value = 2582
$ws.Cells.Item(136, 13).Value = 578409
